{"js": "const pairs = [\n  [\"2023-03-19 Sunday\", \"2023-03-20 Monday\"],\n  [\"81\u00d796=7776\", \"43\u00d761=2623\"],\n  [\"36\u00d755=1980\", \"75\u00d758=4350\"],\n  [\"24\u00d789=2136\", \"43\u00d743=1849\"],\n  [\"39\u00d746=1794\", \"35\u00d755=1925\"],\n  [\"43\u00d733=1419\", \"36\u00d781=2916\"],\n  [\"27\u00d721=567\", \"94\u00d719=1786\"],\n  [\"21\u00d728=588\", \"15\u00d729=435\"],\n  [\"78\u00d794=7332\", \"61\u00d773=4453\"],\n  [\"57\u00d729=1653\", \"60\u00d791=5460\"],\n  [\"53\u00d754=2862\", \"47\u00d788=4136\"],\n  [\"66\u00d786=5676\", \"29\u00d737=1073\"],\n  [\"56\u00d752=2912\", \"22\u00d741=902\"],\n  [\"19\u00d728=532\", \"79\u00d766=5214\"],\n  [\"84\u00d732=2688\", \"50\u00d747=2350\"],\n  [\"22\u00d747=1034\", \"11\u00d758=638\"],\n  [\"94\u00d753=4982\", \"77\u00d788=6776\"],\n  [\"45\u00d728=1260\", \"66\u00d790=5940\"],\n  [\"23\u00d736=828\", \"60\u00d721=1260\"],\n  [\"41\u00d798=4018\", \"54\u00d756=3024\"],\n  [\"81\u00d719=1539\", \"90\u00d773=6570\"],\n  [\"48\u00d730=1440\", \"35\u00d710=350\"],\n  [\"26\u00d790=2340\", \"38\u00d757=2166\"],\n  [\"12\u00d730=360\", \"31\u00d739=1209\"],\n  [\"65\u00d789=5785\", \"49\u00d743=2107\"],\n  [\"90\u00d742=3780\", \"22\u00d763=1386\"],\n  [\"23\u00d747=1081\", \"100\u00d727=2700\"],\n  [\"21\u00d796=2016\", \"67\u00d791=6097\"],\n  [\"84\u00d785=7140\", \"78\u00d750=3900\"],\n  [\"65\u00d791=5915\", \"11\u00d740=440\"],\n  [\"30\u00d766=1980\", \"81\u00d798=7938\"],\n  [\"19\u00d750=950\", \"48\u00d758=2784\"],\n  [\"67\u00d726=1742\", \"19\u00d799=1881\"],\n  [\"75\u00d717=1275\", \"54\u00d788=4752\"],\n  [\"20\u00d757=1140\", \"14\u00d724=336\"],\n  [\"36\u00d752=1872\", \"77\u00d739=3003\"],\n  [\"42\u00d764=2688\", \"10\u00d729=290\"],\n  [\"100\u00d798=9800\", \"45\u00d743=1935\"],\n  [\"59\u00d742=2478\", \"77\u00d764=4928\"],\n  [\"92\u00d794=8648\", \"63\u00d785=5355\"],\n  [\"87\u00d766=5742\", \"28\u00d728=784\"],\n  [\"33\u00d767=2211\", \"54\u00d760=3240\"],\n  [\"30\u00d794=2820\", \"36\u00d739=1404\"],\n  [\"53\u00d761=3233\", \"36\u00d773=2628\"],\n  [\"64\u00d777=4928\", \"72\u00d785=6120\"],\n  [\"40\u00d760=2400\", \"92\u00d746=4232\"],\n  [\"44\u00d710=440\", \"15\u00d765=975\"],\n  [\"58\u00d754=3132\", \"31\u00d749=1519\"],\n  [\"58\u00d758=3364\", \"61\u00d792=5612\"],\n  [\"38\u00d714=532\", \"66\u00d740=2640\"],\n  [\"32\u00d756=1792\", \"63\u00d762=3906\"],\n  [\"84\u00d745=3780\", \"37\u00d785=3145\"],\n  [\"59\u00d723=1357\", \"45\u00d747=2115\"],\n  [\"81\u00d741=3321\", \"90\u00d736=3240\"],\n  [\"53\u00d712=636\", \"72\u00d745=3240\"],\n  [\"47\u00d794=4418\", \"85\u00d756=4760\"],\n  [\"70\u00d745=3150\", \"26\u00d733=858\"],\n  [\"13\u00d764=832\", \"52\u00d770=3640\"],\n  [\"23\u00d754=1242\", \"77\u00d792=7084\"],\n  [\"30\u00d721=630\", \"33\u00d795=3135\"],\n  [\"60\u00d711=660\", \"34\u00d766=2244\"],\n  [\"97\u00d797=9409\", \"68\u00d771=4828\"],\n  [\"78\u00d760=4680\", \"98\u00d753=5194\"],\n  [\"73\u00d783=6059\", \"65\u00d750=3250\"],\n  [\"40\u00d753=2120\", \"36\u00d762=2232\"],\n  [\"17\u00d773=1241\", \"64\u00d776=4864\"],\n  [\"90\u00d735=3150\", \"63\u00d730=1890\"],\n  [\"11\u00d794=1034\", \"51\u00d775=3825\"],\n  [\"67\u00d739=2613\", \"16\u00d799=1584\"],\n  [\"87\u00d770=6090\", \"76\u00d759=4484\"],\n  [\"96\u00d736=3456\", \"24\u00d763=1512\"],\n  [\"71\u00d770=4970\", \"18\u00d788=1584\"],\n  [\"96\u00d786=8256\", \"57\u00d759=3363\"],\n  [\"29\u00d773=2117\", \"86\u00d793=7998\"],\n  [\"43\u00d730=1290\", \"96\u00d743=4128\"],\n  [\"30\u00d714=420\", \"21\u00d711=231\"],\n  [\"14\u00d738=532\", \"36\u00d780=2880\"],\n  [\"37\u00d747=1739\", \"12\u00d726=312\"],\n  [\"55\u00d767=3685\", \"44\u00d775=3300\"],\n  [\"20\u00d727=540\", \"17\u00d745=765\"],\n  [\"18\u00d713=234\", \"29\u00d758=1682\"],\n  [\"24\u00d779=1896\", \"84\u00d779=6636\"],\n  [\"11\u00d751=561\", \"17\u00d736=612\"],\n  [\"17\u00d772=1224\", \"37\u00d713=481\"],\n  [\"18\u00d722=396\", \"78\u00d730=2340\"],\n  [\"33\u00d789=2937\", \"17\u00d725=425\"],\n  [\"23\u00d787=2001\", \"71\u00d727=1917\"],\n  [\"15\u00d787=1305\", \"98\u00d798=9604\"],\n  [\"81\u00d752=4212\", \"60\u00d719=1140\"],\n  [\"95\u00d735=3325\", \"47\u00d726=1222\"],\n  [\"90\u00d732=2880\", \"21\u00d779=1659\"],\n  [\"74\u00d789=6586\", \"59\u00d760=3540\"],\n  [\"10\u00d758=580\", \"83\u00d751=4233\"],\n  [\"38\u00d759=2242\", \"65\u00d755=3575\"],\n  [\"63\u00d731=1953\", \"37\u00d722=814\"],\n  [\"78\u00d784=6552\", \"59\u00d726=1534\"],\n  [\"66\u00d795=6270\", \"54\u00d759=3186\"],\n  [\"61\u00d742=2562\", \"74\u00d775=5550\"],\n  [\"93\u00d739=3627\", \"95\u00d710=950\"],\n  [\"68\u00d765=4420\", \"63\u00d745=2835\"],\n  [\"34\u00d798=3332\", \"29\u00d720=580\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, 'Replace');\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2023-03-19 Sunday\", \"2023-03-20 Monday\")\n    ,@(\"81\u00d796=7776\", \"43\u00d761=2623\")\n    ,@(\"36\u00d755=1980\", \"75\u00d758=4350\")\n    ,@(\"24\u00d789=2136\", \"43\u00d743=1849\")\n    ,@(\"39\u00d746=1794\", \"35\u00d755=1925\")\n    ,@(\"43\u00d733=1419\", \"36\u00d781=2916\")\n    ,@(\"27\u00d721=567\", \"94\u00d719=1786\")\n    ,@(\"21\u00d728=588\", \"15\u00d729=435\")\n    ,@(\"78\u00d794=7332\", \"61\u00d773=4453\")\n    ,@(\"57\u00d729=1653\", \"60\u00d791=5460\")\n    ,@(\"53\u00d754=2862\", \"47\u00d788=4136\")\n    ,@(\"66\u00d786=5676\", \"29\u00d737=1073\")\n    ,@(\"56\u00d752=2912\", \"22\u00d741=902\")\n    ,@(\"19\u00d728=532\", \"79\u00d766=5214\")\n    ,@(\"84\u00d732=2688\", \"50\u00d747=2350\")\n    ,@(\"22\u00d747=1034\", \"11\u00d758=638\")\n    ,@(\"94\u00d753=4982\", \"77\u00d788=6776\")\n    ,@(\"45\u00d728=1260\", \"66\u00d790=5940\")\n    ,@(\"23\u00d736=828\", \"60\u00d721=1260\")\n    ,@(\"41\u00d798=4018\", \"54\u00d756=3024\")\n    ,@(\"81\u00d719=1539\", \"90\u00d773=6570\")\n    ,@(\"48\u00d730=1440\", \"35\u00d710=350\")\n    ,@(\"26\u00d790=2340\", \"38\u00d757=2166\")\n    ,@(\"12\u00d730=360\", \"31\u00d739=1209\")\n    ,@(\"65\u00d789=5785\", \"49\u00d743=2107\")\n    ,@(\"90\u00d742=3780\", \"22\u00d763=1386\")\n    ,@(\"23\u00d747=1081\", \"100\u00d727=2700\")\n    ,@(\"21\u00d796=2016\", \"67\u00d791=6097\")\n    ,@(\"84\u00d785=7140\", \"78\u00d750=3900\")\n    ,@(\"65\u00d791=5915\", \"11\u00d740=440\")\n    ,@(\"30\u00d766=1980\", \"81\u00d798=7938\")\n    ,@(\"19\u00d750=950\", \"48\u00d758=2784\")\n    ,@(\"67\u00d726=1742\", \"19\u00d799=1881\")\n    ,@(\"75\u00d717=1275\", \"54\u00d788=4752\")\n    ,@(\"20\u00d757=1140\", \"14\u00d724=336\")\n    ,@(\"36\u00d752=1872\", \"77\u00d739=3003\")\n    ,@(\"42\u00d764=2688\", \"10\u00d729=290\")\n    ,@(\"100\u00d798=9800\", \"45\u00d743=1935\")\n    ,@(\"59\u00d742=2478\", \"77\u00d764=4928\")\n    ,@(\"92\u00d794=8648\", \"63\u00d785=5355\")\n    ,@(\"87\u00d766=5742\", \"28\u00d728=784\")\n    ,@(\"33\u00d767=2211\", \"54\u00d760=3240\")\n    ,@(\"30\u00d794=2820\", \"36\u00d739=1404\")\n    ,@(\"53\u00d761=3233\", \"36\u00d773=2628\")\n    ,@(\"64\u00d777=4928\", \"72\u00d785=6120\")\n    ,@(\"40\u00d760=2400\", \"92\u00d746=4232\")\n    ,@(\"44\u00d710=440\", \"15\u00d765=975\")\n    ,@(\"58\u00d754=3132\", \"31\u00d749=1519\")\n    ,@(\"58\u00d758=3364\", \"61\u00d792=5612\")\n    ,@(\"38\u00d714=532\", \"66\u00d740=2640\")\n    ,@(\"32\u00d756=1792\", \"63\u00d762=3906\")\n    ,@(\"84\u00d745=3780\", \"37\u00d785=3145\")\n    ,@(\"59\u00d723=1357\", \"45\u00d747=2115\")\n    ,@(\"81\u00d741=3321\", \"90\u00d736=3240\")\n    ,@(\"53\u00d712=636\", \"72\u00d745=3240\")\n    ,@(\"47\u00d794=4418\", \"85\u00d756=4760\")\n    ,@(\"70\u00d745=3150\", \"26\u00d733=858\")\n    ,@(\"13\u00d764=832\", \"52\u00d770=3640\")\n    ,@(\"23\u00d754=1242\", \"77\u00d792=7084\")\n    ,@(\"30\u00d721=630\", \"33\u00d795=3135\")\n    ,@(\"60\u00d711=660\", \"34\u00d766=2244\")\n    ,@(\"97\u00d797=9409\", \"68\u00d771=4828\")\n    ,@(\"78\u00d760=4680\", \"98\u00d753=5194\")\n    ,@(\"73\u00d783=6059\", \"65\u00d750=3250\")\n    ,@(\"40\u00d753=2120\", \"36\u00d762=2232\")\n    ,@(\"17\u00d773=1241\", \"64\u00d776=4864\")\n    ,@(\"90\u00d735=3150\", \"63\u00d730=1890\")\n    ,@(\"11\u00d794=1034\", \"51\u00d775=3825\")\n    ,@(\"67\u00d739=2613\", \"16\u00d799=1584\")\n    ,@(\"87\u00d770=6090\", \"76\u00d759=4484\")\n    ,@(\"96\u00d736=3456\", \"24\u00d763=1512\")\n    ,@(\"71\u00d770=4970\", \"18\u00d788=1584\")\n    ,@(\"96\u00d786=8256\", \"57\u00d759=3363\")\n    ,@(\"29\u00d773=2117\", \"86\u00d793=7998\")\n    ,@(\"43\u00d730=1290\", \"96\u00d743=4128\")\n    ,@(\"30\u00d714=420\", \"21\u00d711=231\")\n    ,@(\"14\u00d738=532\", \"36\u00d780=2880\")\n    ,@(\"37\u00d747=1739\", \"12\u00d726=312\")\n    ,@(\"55\u00d767=3685\", \"44\u00d775=3300\")\n    ,@(\"20\u00d727=540\", \"17\u00d745=765\")\n    ,@(\"18\u00d713=234\", \"29\u00d758=1682\")\n    ,@(\"24\u00d779=1896\", \"84\u00d779=6636\")\n    ,@(\"11\u00d751=561\", \"17\u00d736=612\")\n    ,@(\"17\u00d772=1224\", \"37\u00d713=481\")\n    ,@(\"18\u00d722=396\", \"78\u00d730=2340\")\n    ,@(\"33\u00d789=2937\", \"17\u00d725=425\")\n    ,@(\"23\u00d787=2001\", \"71\u00d727=1917\")\n    ,@(\"15\u00d787=1305\", \"98\u00d798=9604\")\n    ,@(\"81\u00d752=4212\", \"60\u00d719=1140\")\n    ,@(\"95\u00d735=3325\", \"47\u00d726=1222\")\n    ,@(\"90\u00d732=2880\", \"21\u00d779=1659\")\n    ,@(\"74\u00d789=6586\", \"59\u00d760=3540\")\n    ,@(\"10\u00d758=580\", \"83\u00d751=4233\")\n    ,@(\"38\u00d759=2242\", \"65\u00d755=3575\")\n    ,@(\"63\u00d731=1953\", \"37\u00d722=814\")\n    ,@(\"78\u00d784=6552\", \"59\u00d726=1534\")\n    ,@(\"66\u00d795=6270\", \"54\u00d759=3186\")\n    ,@(\"61\u00d742=2562\", \"74\u00d775=5550\")\n    ,@(\"93\u00d739=3627\", \"95\u00d710=950\")\n    ,@(\"68\u00d765=4420\", \"63\u00d745=2835\")\n    ,@(\"34\u00d798=3332\", \"29\u00d720=580\")\n)\n\n$notFound = @()\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        $notFound += $oldText\n    }\n}\n\nif ($notFound.Count -gt 0) {\n    Write-Output (\"NOT FOUND: \" + ($notFound -join \", \"))\n} else {\n    Write-Output \"All replacements applied.\"\n}"}
